$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Paragraph 32 (empty "List Paragraph" bullet at ilvl=1/numId=2) loses
#    its list numbering / ListParagraph style, becoming a plain paragraph
#    that still carries the <w:rPr><w:lang w:val="vi-VN"/></w:rPr> mark.
# ---------------------------------------------------------------------
$pEmpty = $d.Paragraphs(32)
$pEmpty.Style = "Normal"
$pEmpty.Range.LanguageID = "vi-VN"

# ---------------------------------------------------------------------
# 2) Paragraph 45 ("console.log() nằm trong Inspect trang > console")
#    currently split across three runs; merge the first two into a
#    single run while leaving the trailing " > console" run untouched
#    (and unmerged). We briefly give that trailing run a distinguishing
#    format (Bold) so the engine's same-format run consolidation can't
#    fold it into the merge, then revert the formatting afterwards.
# ---------------------------------------------------------------------
$p45 = $d.Paragraphs(45)
$start45 = $p45.Range.Start

$full0 = $p45.Range.Text
$tail = " > console"
$idxTail = $full0.IndexOf($tail)
$lenTail = $tail.Length
$rTail = $d.Range($start45 + $idxTail, $start45 + $idxTail + $lenTail)
$rTail.Bold = 1

$full1 = $p45.Range.Text
$head = "console.log()"
$mid = "Inspect trang"
$idxHead = $full1.IndexOf($head)
$idxMidEnd = $full1.IndexOf($mid) + $mid.Length
$rMerge = $d.Range($start45 + $idxHead, $start45 + $idxMidEnd)
# force a genuine content change (no-op text-only assignment does not
# trigger run consolidation) ...
$rMerge.Text = "console.log() nằm trong Inspect trang!"

# ... then set it back to the real merged text.
$full2 = $p45.Range.Text
$idxHead2 = $full2.IndexOf($head)
$idxMidEnd2 = $full2.IndexOf("$mid!") + "$mid!".Length
$rMerge2 = $d.Range($start45 + $idxHead2, $start45 + $idxMidEnd2)
$rMerge2.Text = "console.log() nằm trong Inspect trang"

# restore the trailing run's original formatting
$full3 = $p45.Range.Text
$idxTail2 = $full3.IndexOf($tail)
$rTail2 = $d.Range($start45 + $idxTail2, $start45 + $idxTail2 + $lenTail)
$rTail2.Bold = 0

# ---------------------------------------------------------------------
# 3) Paragraph 52 (last paragraph) "document.getElementById(\u201cid\u201d)."
#    + "value áp dụng ..." split across two runs; merge into one run.
#    This paragraph's two runs are the entirety of its text, so a
#    straightforward replace is safe once we force a real content
#    change (same no-op caveat as above).
# ---------------------------------------------------------------------
$p52 = $d.Paragraphs($d.Paragraphs.Count)
$start52 = $p52.Range.Start
$end52 = $p52.Range.End
$target = "document.getElementById(“id”).value áp dụng cho thẻ input tự tạo, lấy giá trị của input thông qua id"

$rFull = $d.Range($start52, $end52 - 1)
$rFull.Text = $target + "!"

$p52b = $d.Paragraphs($d.Paragraphs.Count)
$rFull2 = $d.Range($p52b.Range.Start, $p52b.Range.End - 1)
$rFull2.Text = $target

Write-Output "edits applied"
